$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update conditions for running the script (new shared strings created first)
$ws.Range("D5").Value = "Every 24 hours, at GMT-10"
$ws.Range("D6").Value = "Every 24 hours, at GMT-10"
$ws.Range("D9").Value = "Right after running the train_users script at GMT-10"

# Update script names
$ws.Range("B5").Value = "train_visitors"
$ws.Range("B6").Value = "train_users"

# Update selection to match final state
$ws.Range("B15").Select()
